$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 214 (shifts existing rows 214-320 down to 217-323)
$ws.Rows(214).Insert()
$ws.Rows(214).Insert()
$ws.Rows(214).Insert()

# Row 214
$ws.Cells.Item(214, 1).Value = 5
$ws.Cells.Item(214, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(214, 3).Value = 'Maule'
$ws.Cells.Item(214, 4).Value = 44813
$ws.Cells.Item(214, 5).Value = 7
$ws.Cells.Item(214, 6).Value = 100112045
$ws.Cells.Item(214, 7).Value = 'Zapallo'
$ws.Cells.Item(214, 8).Value = 'Camote'
$ws.Cells.Item(214, 9).Value = '1a (guarda)'
$ws.Cells.Item(214, 10).Value = 700
$ws.Cells.Item(214, 11).Value = 1000
$ws.Cells.Item(214, 12).Value = 1000
$ws.Cells.Item(214, 13).Value = 1000
$ws.Cells.Item(214, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(214, 15).Value = 'Región del Maule'
$ws.Cells.Item(214, 16).Value = 1000
$ws.Cells.Item(214, 17).Value = 1
$ws.Cells.Item(214, 18).Value = 'Hortaliza'

# Row 215
$ws.Cells.Item(215, 1).Value = 5
$ws.Cells.Item(215, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(215, 3).Value = 'Maule'
$ws.Cells.Item(215, 4).Value = 44813
$ws.Cells.Item(215, 5).Value = 7
$ws.Cells.Item(215, 6).Value = 100112045
$ws.Cells.Item(215, 7).Value = 'Zapallo'
$ws.Cells.Item(215, 8).Value = 'Camote'
$ws.Cells.Item(215, 9).Value = '2a (guarda)'
$ws.Cells.Item(215, 10).Value = 90
$ws.Cells.Item(215, 11).Value = 700
$ws.Cells.Item(215, 12).Value = 700
$ws.Cells.Item(215, 13).Value = 700
$ws.Cells.Item(215, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(215, 15).Value = 'Región del Maule'
$ws.Cells.Item(215, 16).Value = 700
$ws.Cells.Item(215, 17).Value = 1
$ws.Cells.Item(215, 18).Value = 'Hortaliza'

# Row 216
$ws.Cells.Item(216, 1).Value = 5
$ws.Cells.Item(216, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(216, 3).Value = 'Maule'
$ws.Cells.Item(216, 4).Value = 44813
$ws.Cells.Item(216, 5).Value = 7
$ws.Cells.Item(216, 6).Value = 100112045
$ws.Cells.Item(216, 7).Value = 'Zapallo'
$ws.Cells.Item(216, 8).Value = 'Paine'
$ws.Cells.Item(216, 9).Value = '1a (guarda)'
$ws.Cells.Item(216, 10).Value = 2000
$ws.Cells.Item(216, 11).Value = 250
$ws.Cells.Item(216, 12).Value = 250
$ws.Cells.Item(216, 13).Value = 250
$ws.Cells.Item(216, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(216, 15).Value = 'Región del Maule'
$ws.Cells.Item(216, 16).Value = 250
$ws.Cells.Item(216, 17).Value = 1
$ws.Cells.Item(216, 18).Value = 'Hortaliza'
